# Automatische test-sync: 2025-08-05 19:51:50
# Adds a new test-mail row (row 55) to the "Logs" sheet, extends the
# conditional-formatting ranges that tracked the used range, and bumps the
# matching category counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 55 ------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(55, 1).Value = "Heb je ergens de CE-certificaten van dit product?"
$logs.Cells.Item(55, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(55, 3).Value = "Testmail #14: Heb je ergens de CE-certificaten van dit product?"
$logs.Cells.Item(55, 4).Value = "Kwaliteit / Certificaten"
$logs.Cells.Item(55, 5).Value = "Bedankt, we hebben dit doorgestuurd naar kwaliteit@bedrijf.nl."
$logs.Cells.Item(55, 6).Value = "2025-08-05 19:51:36"
$logs.Cells.Item(55, 7).Value = "Ja"
$logs.Cells.Item(55, 8).Value = "Ja"
$logs.Cells.Item(55, 9).Value = "Nee"
$logs.Cells.Item(55, 10).Value = "Nee"

# --- extend the conditional formatting ranges from row 54 to row 55 ----------
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "54")
    $newRange = $logs.Range($col + "2:" + $col + "55")
    $cfs = $oldRange.FormatConditions
    for ($i = 1; $i -le $cfs.Count; $i++) {
        $cfs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: bump the "Kwaliteit / Certificaten" count --------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(10, 2).Value = 2
